$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.068707790979985361
$ws.Cells.Item(1, 2).Value = 0.068650238769620842
$ws.Cells.Item(2, 1).Value = 0.015387258832618045
$ws.Cells.Item(2, 2).Value = -0.015601913886255758
$ws.Cells.Item(3, 1).Value = 0.11853066833615244
$ws.Cells.Item(3, 2).Value = -0.11887461911758379
$ws.Cells.Item(4, 1).Value = -0.16511725079794104
$ws.Cells.Item(4, 2).Value = 0.16445802748173932
$ws.Cells.Item(5, 1).Value = -0.15845802795394093
$ws.Cells.Item(5, 2).Value = 0.15714402695076313
$ws.Cells.Item(6, 1).Value = -0.058978007070919158
$ws.Cells.Item(6, 2).Value = 0.058935100551482655
$ws.Cells.Item(7, 1).Value = -0.038935101119040638
$ws.Cells.Item(7, 2).Value = 0.038866935294416649
$ws.Cells.Item(8, 1).Value = -0.018866935864645384
$ws.Cells.Item(8, 2).Value = 0.018837468398170465
$ws.Cells.Item(9, 1).Value = -0.012837468890983139
$ws.Cells.Item(9, 2).Value = 0.012817509510648861
$ws.Cells.Item(10, 1).Value = -0.0068175100047227488
$ws.Cells.Item(10, 2).Value = 0.0068179065507933956
$ws.Cells.Item(11, 1).Value = -0.0023179070364243159
$ws.Cells.Item(11, 2).Value = 0.0023134417709620436
$ws.Cells.Item(12, 1).Value = 0.0036865577347691136
$ws.Cells.Item(12, 2).Value = -0.0037263154920341712
$ws.Cells.Item(13, 1).Value = 0.0097263149982920183
$ws.Cells.Item(13, 2).Value = -0.0097488374661702792
$ws.Cells.Item(14, 1).Value = 0.021748836938899174
$ws.Cells.Item(14, 2).Value = -0.021804246266353466
$ws.Cells.Item(15, 1).Value = -0.021049532314417974
$ws.Cells.Item(15, 2).Value = 0.021025874001995959
$ws.Cells.Item(16, 1).Value = -0.015025874495488978
$ws.Cells.Item(16, 2).Value = 0.015004123167893457
$ws.Cells.Item(17, 1).Value = -0.0090041236631996924
$ws.Cells.Item(17, 2).Value = 0.0089999994873064537
$ws.Cells.Item(18, 1).Value = -0.037975563363850995
$ws.Cells.Item(18, 2).Value = 0.037961044601715344
$ws.Cells.Item(19, 1).Value = -0.027096470583364507
$ws.Cells.Item(19, 2).Value = 0.027013615788022083
$ws.Cells.Item(20, 1).Value = -0.018013616279629829
$ws.Cells.Item(20, 2).Value = 0.018004259459363681
$ws.Cells.Item(21, 1).Value = -0.0090042599515474109
$ws.Cells.Item(21, 2).Value = 0.0089999995074254713
$ws.Cells.Item(22, 1).Value = -0.093928754936747794
$ws.Cells.Item(22, 2).Value = 0.093621972661406261
$ws.Cells.Item(23, 1).Value = -0.084621973152301244
$ws.Cells.Item(23, 2).Value = 0.084124289635691518
$ws.Cells.Item(24, 1).Value = -0.042124290321613245
$ws.Cells.Item(24, 2).Value = 0.041999999310557534
$ws.Cells.Item(25, 1).Value = -0.092166020366956758
$ws.Cells.Item(25, 2).Value = 0.09205889321120253
$ws.Cells.Item(26, 1).Value = -0.08605889369983899
$ws.Cells.Item(26, 2).Value = 0.085926245681783087
$ws.Cells.Item(27, 1).Value = -0.079926246172705717
$ws.Cells.Item(27, 2).Value = 0.079490865523421661
$ws.Cells.Item(28, 1).Value = -0.073490866023987245
$ws.Cells.Item(28, 2).Value = 0.073208677995517668
$ws.Cells.Item(29, 1).Value = -0.061208678536505801
$ws.Cells.Item(29, 2).Value = 0.061083527011989602
$ws.Cells.Item(30, 1).Value = -0.041083527601204928
$ws.Cells.Item(30, 2).Value = 0.040683502020741003
$ws.Cells.Item(31, 1).Value = -0.027018554516519799
$ws.Cells.Item(31, 2).Value = 0.027000781369862281
$ws.Cells.Item(32, 1).Value = -0.0060007819699121612
$ws.Cells.Item(32, 2).Value = 0.0059999994848576321

Write-Output "Done updating A1:B32"
